# Auto-generated edit script: applies the numeric cell updates
# described by the upstream OOXML diff for Sheets/Maduin_Profits.xlsx.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 9916.333000000001
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("H132").Value = 2479
$ws.Range("I132").Value = 2562.2727
$ws.Range("K132").Value = 7686.8181
$ws.Range("M132").Value = -5156.8181
$ws.Range("H137").Value = 2313.4
$ws.Range("I137").Value = 1877.7142
$ws.Range("J137").Value = 3330
$ws.Range("K137").Value = 5633.142599999999
$ws.Range("L137").Value = 9990
$ws.Range("M137").Value = -3083.142599999999
$ws.Range("N137").Value = -15090
$ws.Range("H138").Value = 13437.533
$ws.Range("J138").Value = 13437.533
$ws.Range("L138").Value = 40312.599
$ws.Range("N138").Value = -50592.599

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11642.667
$ws.Range("I32").Value = 10640.037
$ws.Range("K32").Value = 10640.037
$ws.Range("M32").Value = -10353.037
$ws.Range("H61").Value = 5627.8335
$ws.Range("I61").Value = 4384.5
$ws.Range("J61").Value = 6249.5
$ws.Range("K61").Value = 4384.5
$ws.Range("L61").Value = 6249.5
$ws.Range("M61").Value = -4172.5
$ws.Range("N61").Value = -6673.5
$ws.Range("H74").Value = 4125
$ws.Range("I74").Value = 4125
$ws.Range("K74").Value = 4125
$ws.Range("M74").Value = -3251
$ws.Range("H77").Value = 4125
$ws.Range("I77").Value = 4125
$ws.Range("K77").Value = 20625
$ws.Range("M77").Value = -16257
$ws.Range("H136").Value = 5627.8335
$ws.Range("I136").Value = 4384.5
$ws.Range("J136").Value = 6249.5
$ws.Range("K136").Value = 13153.5
$ws.Range("L136").Value = 18748.5
$ws.Range("M136").Value = -10603.5
$ws.Range("N136").Value = -23848.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 98999
$ws.Range("J59").Value = 98999
$ws.Range("L59").Value = 98999
$ws.Range("N59").Value = -100693
$ws.Range("H61").Value = 48000
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 48000
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 48000
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -48626
$ws.Range("H94").Value = 1424.619
$ws.Range("I94").Value = 1145.85
$ws.Range("K94").Value = 1145.85
$ws.Range("M94").Value = -694.8499999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3364.5454
$ws.Range("I31").Value = 3168
$ws.Range("K31").Value = 3168
$ws.Range("M31").Value = -2873
$ws.Range("H34").Value = 3364.5454
$ws.Range("I34").Value = 3168
$ws.Range("K34").Value = 3168
$ws.Range("M34").Value = -2966
$ws.Range("H45").Value = 26500
$ws.Range("I45").Value = 7000
$ws.Range("J45").Value = 46000
$ws.Range("K45").Value = 7000
$ws.Range("L45").Value = 46000
$ws.Range("M45").Value = -6407
$ws.Range("N45").Value = -47186
$ws.Range("H106").Value = 500000
$ws.Range("J106").Value = 500000
$ws.Range("L106").Value = 500000
$ws.Range("N106").Value = -502524
$ws.Range("H132").Value = 4014.7693
$ws.Range("I132").Value = 3099.3333
$ws.Range("K132").Value = 9297.999899999999
$ws.Range("M132").Value = -6767.999899999999
$ws.Range("H134").Value = 3926
$ws.Range("I134").Value = 3462.842
$ws.Range("J134").Value = 5183.143
$ws.Range("K134").Value = 10388.526
$ws.Range("L134").Value = 15549.429
$ws.Range("M134").Value = -7853.526
$ws.Range("N134").Value = -20619.429

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 520.3
$ws.Range("I5").Value = 467
$ws.Range("J5").Value = 1000
$ws.Range("K5").Value = 1401
$ws.Range("L5").Value = 3000
$ws.Range("M5").Value = -1289
$ws.Range("N5").Value = -3224
$ws.Range("H6").Value = 168.33333
$ws.Range("J6").Value = 3
$ws.Range("L6").Value = 9
$ws.Range("N6").Value = -235
$ws.Range("H68").Value = 1600
$ws.Range("I68").Value = 1600
$ws.Range("K68").Value = 4800
$ws.Range("M68").Value = -3989
$ws.Range("H71").Value = 1600
$ws.Range("I71").Value = 1600
$ws.Range("K71").Value = 14400
$ws.Range("M71").Value = -10344
$ws.Range("H113").Value = 850
$ws.Range("I113").Value = 850
$ws.Range("K113").Value = 2550
$ws.Range("M113").Value = -380
$ws.Range("H132").Value = 3171.6086
$ws.Range("I132").Value = 1998
$ws.Range("J132").Value = 3224.9546
$ws.Range("K132").Value = 17982
$ws.Range("L132").Value = 29024.5914
$ws.Range("M132").Value = -15452
$ws.Range("N132").Value = -34084.5914
$ws.Range("H135").Value = 520.3
$ws.Range("I135").Value = 467
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 4203
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = -1668
$ws.Range("N135").Value = -14070

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 5000
$ws.Range("J40").Value = 5000
$ws.Range("L40").Value = 5000
$ws.Range("N40").Value = -5302
$ws.Range("H44").Value = 15000
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").ClearContents()
$ws.Range("H57").Value = 15249.75
$ws.Range("I57").Value = 3000
$ws.Range("K57").Value = 3000
$ws.Range("M57").Value = -2180
$ws.Range("H102").Value = 1766
$ws.Range("J102").Value = 2999
$ws.Range("L102").Value = 2999
$ws.Range("N102").Value = -6243
$ws.Range("H107").Value = 875
$ws.Range("I107").Value = 875
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 875
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1045
$ws.Range("N107").ClearContents()
$ws.Range("H132").Value = 3350.2632
$ws.Range("I132").Value = 2811.25
$ws.Range("K132").Value = 8433.75
$ws.Range("M132").Value = -5903.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 988.8889
$ws.Range("I22").Value = 988.8889
$ws.Range("K22").Value = 988.8889
$ws.Range("M22").Value = -693.8889
$ws.Range("H27").Value = 988.8889
$ws.Range("I27").Value = 988.8889
$ws.Range("K27").Value = 988.8889
$ws.Range("M27").Value = -881.8889
$ws.Range("H46").Value = 1442
$ws.Range("I46").Value = 1296.5
$ws.Range("K46").Value = 1296.5
$ws.Range("M46").Value = -1108.5
$ws.Range("I122").Value = 8564.714
$ws.Range("J122").Value = 8000
$ws.Range("K122").Value = 25694.142
$ws.Range("L122").Value = 24000
$ws.Range("M122").Value = -23244.142
$ws.Range("N122").Value = -28900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 75995
$ws.Range("J64").Value = 75995
$ws.Range("L64").Value = 75995
$ws.Range("N64").Value = -76491
$ws.Range("H67").Value = 75995
$ws.Range("J67").Value = 75995
$ws.Range("L67").Value = 75995
$ws.Range("N67").Value = -77711
$ws.Range("H104").Value = 18500
$ws.Range("J104").Value = 18500
$ws.Range("L104").Value = 18500
$ws.Range("N104").Value = -25488
$ws.Range("H132").Value = 1939.3043
$ws.Range("I132").Value = 1106.5555
$ws.Range("K132").Value = 3319.6665
$ws.Range("M132").Value = -789.6664999999998
$ws.Range("H136").Value = 66274.625
$ws.Range("I136").Value = 102039.6
$ws.Range("J136").Value = 6666.3335
$ws.Range("K136").Value = 306118.8
$ws.Range("L136").Value = 19999.0005
$ws.Range("M136").Value = -303568.8
$ws.Range("N136").Value = -25099.0005
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

